# Update cars.xlsx: refresh listing data (newer data / PMGWeb upload)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Nissan NP300 listing replaced with a newer Toyota Hilux listing
$ws.Range("A2").Value = "Toyota Hilux 2.8GD-6 Double Cab 4x4 Legend 50 Auto"
$ws.Range("B2").Value = 10141375
$ws.Range("C2").Value = "UG5438"
$ws.Range("D2").Value = "Used"
$ws.Range("E2").Value = 2020
$ws.Range("F2").Value = "75 000 Km"
$ws.Range("G2").Value = "R 619 900"

# Remove two rows that no longer belong in the feed:
# row 18 "Ford Ranger 3.2 TDCi XLT 4x4 Auto Double-Cab" and
# row 15 "Toyota Fortuner 2.5 D-4D Raised Body Auto"
# (delete the higher row index first so the lower index is unaffected)
$ws.Rows(18).EntireRow.Delete()
$ws.Rows(15).EntireRow.Delete()
